# Refresh the cryptos price list snapshot (prices + 1h volume deltas).
# D-column values are stored as text (e.g. "28.187.90" is not a valid
# number), so new values are written with a leading quote-prefix escape
# to stop Excel re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'28.226.72"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "`'1.825.43"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").Value = "`'1.000"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "`'339.66"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "`'0.9980"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "`'0.3938"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("D8").Value = "`'0.3507"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").Value = "`'48.20"
$ws.Range("D10").Value = "`'1.202"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Value = "`'0.07605"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "`'0.9980"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "`'22.21"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "`'6.555"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "`'1.826.46"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "`'7.221"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "`'0.06704"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "`'85.64"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "`'0.9973"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "`'17.94"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").Value = "`'6.600"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "`'28.269.42"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").Value = "`'12.79"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "`'2.403"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "`'2.584"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "`'1.497"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "`'21.54"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "`'155.01"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "`'2.033.46"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "`'136.00"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "`'6.240"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "`'4.041"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "`'0.08862"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").Value = "`'13.30"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "`'5.550"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "`'0.02454"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "`'0.6987"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").Value = "`'0.06576"
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("D40").Value = "`'1.616"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "`'0.2232"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "`'8.590"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("D44").Value = "`'14.63"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "`'0.6525"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").Value = "`'132.00"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "`'0.07225"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "`'80.69"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "`'1.251"
$ws.Range("E51").Value = "  +4.50%  "
